$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: the "Horizontal custom rolls" feature is now completed in version 1.5.7.
# Its old requester (Sean Steele, previously in E16) moves to a new row describing
# a follow-up feature request that came out of the horizontal layout change.
$ws.Range("C16").Value = "1.5.7"
$ws.Range("E16").Value = ""

# New row 19: "Edit Saved Rolls" feature request from Sean Steele.
$ws.Range("A19").Value = "Edit Saved Rolls"
$ws.Range("B19").Value = "Allow the user to edit saved rolls easily. Auto move to custom roll tab, and fill out all the die in the roll."
$ws.Range("D19").Value = "Sean Steele - ssteele1812@gmail.com"

# Update selection to reflect where the user left off editing.
$ws.Range("B20").Select()
